# Daily Satellite Data Update
# Adds the 31.12.2025 passes (rows 8 & 9) and refreshes the cloud-cover
# statistics (columns O/P/Q/R) for the existing rows (2-7), including the
# per-cell heat-map fill colors that go with the refreshed numbers, and
# extends the conditional-formatting ranges to cover the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the Oblačnost / Nízká / Střední / Vysoká numbers for the
#    existing rows (2-7) plus their fill colours (heat-map style).
# ---------------------------------------------------------------------

# Colours (converted to the OLE BGR integer Interior.Color expects)
$green      = 7000682    # FF6AD26A - low cloud cover
$blue       = 13143125   # FF558CC8 - baseline (0) Stredni/Vysoka
$purple     = 13473638   # FF6697CD - medium-low cloud cover
$red        = 8351984    # FFF0707F - high cloud cover
$lightblue  = 14927274   # FFAAC5E3
$lighterblue= 15323579   # FFBBD1E9
$medblue    = 14794914   # FFA2C0E1

function Set-Stats($row, $o, $oColor, $p, $pColor, $q, $r) {
    $ws.Range("O$row").Value = $o
    $ws.Range("O$row").Interior.Color = $oColor
    $ws.Range("P$row").Value = $p
    $ws.Range("P$row").Interior.Color = $pColor
    $ws.Range("Q$row").Value = $q
    $ws.Range("Q$row").Interior.Color = $blue
    $ws.Range("R$row").Value = $r
    $ws.Range("R$row").Interior.Color = $blue
}

Set-Stats 2 2  $green      0  $blue       1 0
Set-Stats 3 11 $purple     8  $purple     0 0
Set-Stats 4 57 $red        49 $lightblue  0 0
Set-Stats 5 76 $red        60 $lighterblue 0 0
Set-Stats 6 68 $red        59 $lighterblue 0 0
Set-Stats 7 61 $red        47 $medblue    0 0

# ---------------------------------------------------------------------
# 2) Append the two new passes for 31.12.2025 (rows 8 and 9).
# ---------------------------------------------------------------------

$ws.Range("A8").Value = "31.12.2025"
$ws.Range("B8").Value = 16
$ws.Range("C8").Value = "04:41"
$ws.Range("D8").Value = "10:08"
$ws.Range("E8").Value = "03:56:02"
$ws.Range("F8").Value = "03:58:45"
$ws.Range("G8").Value = "04:01:05"
$ws.Range("H8").Value = "04:03:26"
$ws.Range("I8").Value = "04:06:10"
$ws.Range("J8").Value = "16°"
$ws.Range("K8").Value = "04:01:39"
$ws.Range("L8").Value = -28
$ws.Range("M8").Value = "A+B"
$ws.Range("N8").Value = "2"

$ws.Range("A9").Value = "31.12.2025"
$ws.Range("B9").Value = 31
$ws.Range("C9").Value = "06:31"
$ws.Range("D9").Value = "11:06"
$ws.Range("E9").Value = "05:33:10"
$ws.Range("F9").Value = "05:35:27"
$ws.Range("G9").Value = "05:38:43"
$ws.Range("H9").Value = "05:41:58"
$ws.Range("I9").Value = "05:44:16"
$ws.Range("J9").Value = "10°"
$ws.Range("K9").Value = "05:35:26"
$ws.Range("L9").Value = -13
$ws.Range("M9").Value = "B"
$ws.Range("N9").Value = "3"

# ---------------------------------------------------------------------
# 3) Extend conditional formatting ranges to include the new rows.
# ---------------------------------------------------------------------

$lCond = $ws.Range("L2:L7").FormatConditions
for ($i = 1; $i -le $lCond.Count; $i++) {
    $lCond.Item($i).ModifyAppliesToRange($ws.Range("L2:L9"))
}

$nCond = $ws.Range("N2:N7").FormatConditions
for ($i = 1; $i -le $nCond.Count; $i++) {
    $nCond.Item($i).ModifyAppliesToRange($ws.Range("N2:N9"))
}

Write-Host "Daily satellite data update applied"
